$d = $word.ActiveDocument
$apos = [char]0x2019

# ============================================================
# Edit 2 first (paragraph: "...the only even numbers to pan out in the
# little girl's counting pattern."):
#   change the trailing "pattern" to "method" -> "...counting method."
#   Done before Edit 1 so the pre-existing "_GoBack" bookmark (originally
#   sitting right before the final period) acts as a natural barrier that
#   keeps " method" and "." in separate runs.
# ============================================================

$anchor2 = $d.Content
$anchor2.Find.Execute("pan out in the little girl") | Out-Null
$afterAnchor2 = $anchor2.End

$wordBoundary = $d.Range($afterAnchor2, $d.Content.End)
$wordBoundary.Find.Execute("${apos}s counting") | Out-Null
$splitPos = $wordBoundary.End

# Temporarily bookmark the split point so " pattern" -> " method" becomes
# its own run instead of merging back into "...counting".
$barrier2 = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("ZZ_BARRIER2", $barrier2) | Out-Null

$searchRng2 = $d.Range($splitPos, $d.Content.End)
$searchRng2.Find.Execute(" pattern", $true, $false, $false, $false, $false, `
                          $true, 1, $false, " method", 2) | Out-Null

$d.Bookmarks.Item("ZZ_BARRIER2").Delete()

# ============================================================
# Edit 1 (paragraph: "...I've identified a pattern for each finger; which
# may help toward to a solution."):
#   remove the stray "to " so it reads
#   "...which may help toward a solution."
#   and leave the bookmark "_GoBack" positioned right after "toward".
# ============================================================

# Locate the boundary right after the existing "r; " run so our edit,
# which is scoped to start there, does not reach back into it.
$rAnchor = $d.Content
$rAnchor.Find.Execute("r; ") | Out-Null
$afterR = $rAnchor.End

# Temporarily bookmark that boundary - this stops the engine from
# collapsing the untouched "r; " run into the run we are about to edit.
$barrier1 = $d.Range($afterR, $afterR)
$d.Bookmarks.Add("ZZ_BARRIER1", $barrier1) | Out-Null

$searchRng1 = $d.Range($afterR, $d.Content.End)
$searchRng1.Find.Execute("toward to a solution", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "toward a solution", 2) | Out-Null

$d.Bookmarks.Item("ZZ_BARRIER1").Delete()

# Move "_GoBack" so it sits right after "toward" (before " a solution."),
# matching where the edit actually happened - this is also what removes it
# from its old spot between "method" and "." in the previous paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$gbAnchor = $d.Content
$gbAnchor.Find.Execute("which may help toward") | Out-Null
$gbPos = $gbAnchor.End
$gbRange = $d.Range($gbPos, $gbPos)
$d.Bookmarks.Add("_GoBack", $gbRange) | Out-Null
